$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (shifts existing quarterly data right by 2)
$ws.Range("D:E").EntireColumn.Insert()

# Copy number formats from the (now shifted) adjacent data columns F:G into the
# newly inserted D:E columns so the new cells pick up the same date / number
# formatting as the rest of each row. Done per contiguous data block so that
# label-only rows (5,6,37,79) are not touched.
$ws.Range("F7:G35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:G77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:G102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new columns D (2018-12-31) and E (2018-09-30) with the new
# quarterly figures, and correct a handful of restated prior-quarter values.
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 168500
$ws.Range("E8").Value = 164700
$ws.Range("D9").Value = 47800
$ws.Range("E9").Value = 45000
$ws.Range("D10").Value = 120700
$ws.Range("E10").Value = 119700
$ws.Range("D12").Value = 43400
$ws.Range("E12").Value = 42100
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 270500
$ws.Range("E14").Value = 2900
$ws.Range("D15").Value = 33100
$ws.Range("E15").Value = 42600
$ws.Range("D17").Value = 441900
$ws.Range("E17").Value = 172400
$ws.Range("D18").Value = -273400
$ws.Range("E18").Value = -7700
$ws.Range("D20").Value = -1100
$ws.Range("E20").Value = 1900
$ws.Range("D21").Value = -241400
$ws.Range("E21").Value = 36800
$ws.Range("D22").Value = 12900
$ws.Range("E22").Value = 12400
$ws.Range("D23").Value = -287400
$ws.Range("E23").Value = -18200
$ws.Range("D24").Value = -100
$ws.Range("E24").Value = 3900
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -287400
$ws.Range("E26").Value = -22100
$ws.Range("D27").Value = -287400
$ws.Range("E27").Value = -22100
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = -900
$ws.Range("E29").Value = -800
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 1100
$ws.Range("E32").Value = -1900
$ws.Range("D33").Value = -288200
$ws.Range("E33").Value = -22800
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -288200
$ws.Range("E35").Value = -22800
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 162000
$ws.Range("E41").Value = 149700
$ws.Range("D42").Value = 159000
$ws.Range("E42").Value = 162100
$ws.Range("D43").Value = 154500
$ws.Range("E43").Value = 175300
$ws.Range("D44").Value = 7400
$ws.Range("E44").Value = 8000
$ws.Range("D45").Value = 29200
$ws.Range("E45").Value = 34500
$ws.Range("D46").Value = 512000
$ws.Range("E46").Value = 529500
$ws.Range("D47").Value = 81700
$ws.Range("E47").Value = 78300
$ws.Range("D48").Value = 53600
$ws.Range("E48").Value = 50700
$ws.Range("D49").Value = 2058100
$ws.Range("E49").Value = 2337200
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 54800
$ws.Range("E52").Value = 49100
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 2760300
$ws.Range("E54").Value = 3044800
$ws.Range("D57").Value = 2200
$ws.Range("E57").Value = 10800
$ws.Range("D58").Value = 373400
$ws.Range("E58").Value = 7000
$ws.Range("D59").Value = 148900
$ws.Range("E59").Value = 131400
$ws.Range("D60").Value = 524400
$ws.Range("E60").Value = 149200
$ws.Range("D61").Value = 618800
$ws.Range("E61").Value = 982800
$ws.Range("D62").Value = 124200
$ws.Range("E62").Value = 121200
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 1267400
$ws.Range("E66").Value = 1253200
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -1710600
$ws.Range("E72").Value = -1422400
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 1492900
$ws.Range("E76").Value = 1791600
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -288200
$ws.Range("E81").Value = -22800
$ws.Range("D83").Value = 33100
$ws.Range("E83").Value = 42600
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 44800
$ws.Range("E89").Value = 51900
$ws.Range("D91").Value = -6800
$ws.Range("E91").Value = -2900
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -7700
$ws.Range("E94").Value = -11300
$ws.Range("D96").Value = -22300
$ws.Range("E96").Value = -22300
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -24700
$ws.Range("E100").Value = -23300
$ws.Range("D101").Value = -100
$ws.Range("E101").Value = 100
$ws.Range("D102").Value = 12300
$ws.Range("E102").Value = 17500
$ws.Range("H45").Value = 34800
$ws.Range("H46").Value = 496900
$ws.Range("H52").Value = 71600
$ws.Range("J91").Value = -6000
